$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C (header "Förändrad") holds a date serial value that was bumped
# by one day (45178 -> 45179) for every data row (rows 2 through 338).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45178) {
        $cell.Value2 = 45179
    }
}
